$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number-format/style from row 3 (B and C columns) into the new row 4
# before row 3's own content changes, so row 4 ends up styled like rows 2/3.
$ws.Cells.Item(3, 2).Copy($ws.Cells.Item(4, 2))
$ws.Cells.Item(3, 3).Copy($ws.Cells.Item(4, 3))

# Populate the new 4th data row.
$ws.Cells.Item(4, 1).Value = "one"
$ws.Cells.Item(4, 2).Value = 2
# Set then fix up the bad date-like text value so it reuses/updates the same
# shared string that previously held the malformed date text.
$ws.Cells.Item(4, 3).Value = "201/01/2021"
$ws.Cells.Item(4, 3).Replace("201/01/2021", "202/01/2020")
$ws.Cells.Item(4, 4).Value = "exist"

# The bad date values no longer live in C2/C3 - clear them but keep their
# date number format/style intact.
$ws.Cells.Item(2, 3).ClearContents()
$ws.Cells.Item(3, 3).ClearContents()

# Column C now needs to fit its (longer) textual content.
$ws.Columns.Item(3).ColumnWidth = 11

# Move the active selection down to the new empty row below the data.
$ws.Range("A5").Select()
